$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cell E8 value from "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Set the selection/active cell to E8, matching the diff's sheetView selection
$ws.Range("E8").Select()
